# portfolio.xlsx - "minor adjustment and notes"
# Refresh cached market-data values across watchlist / stocks / portfolio /
# summary sheets (watchlist re-sorts itself by RSI ascending as prices move),
# and drop the now-superseded last row of the summary log.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. watchlist — ticker, price, pct_change, rsi, re-sorted by rsi ascending
# ---------------------------------------------------------------------------
$watchlist = $wb.Worksheets.Item("watchlist")

$rows = @(
    @("WBA",  42.02,               -1.013,   25.58558558558559),
    @("RTX",  59,                   0.2549,  29.25180519280995),
    @("CAT",  108.91,                1.1047, 30.06517427033153),
    @("TRV",  96.84999999999999,    -0.0619, 31.12028158757131),
    @("WMT",  124.73,                0.8327, 31.82879377431907),
    @("VZ",   56.51,                 0.4801, 32.3943661971831),
    @("PG",   116.01,                0.2073, 35.27553527553528),
    @("MCD",  179.24,               -1.4461, 35.68782415411212),
    @("INTC", 58.75,                 1.3106, 38.39483766888486),
    @("DOW",  33.37,                 0.9377, 40.2439024390244),
    @("CSCO", 41.46,                 0.266,  41.97256210604376),
    @("IBM",  122.58,                0.7396, 42.30452674897119),
    @("CVX",  92.89,                 1.5857, 42.91561712846347),
    @("BA",   125.4,                -4.6098, 43.35226946691915),
    @("DIS",  101.06,               -2.0547, 44.59847509654421),
    @("XOM",  44.83,                -0.1114, 44.71853257432005),
    @("NKE",  87.2,                  1.691,  46.7005076142132),
    @("MRK",  78.01000000000001,     1.4962, 47.48294302046838),
    @("GS",   178.3,                 1.2896, 48.03341454925166),
    @("UNH",  293.04,                1.9128, 51.46884021056208),
    @("KO",   45.4,                  0.576,  52.76156264032331),
    @("AXP",  86.73999999999999,     0.4168, 53.31812131583951),
    @("JPM",  92,                   -0.1519, 53.77090564846955),
    @("MSFT", 180.76,                1.0736, 56.60719685515573),
    @("PFE",  38.51,                 2.3658, 59.58188153310105),
    @("AAPL", 297.56,                1.5009, 60.51681550446514),
    @("MMM",  147.43,               -0.6536, 63.22254335260115),
    @("JNJ",  149.5,                 0.8296, 64.44880923152468),
    @("V",    178.44,                1.3,    68.4143455306754),
    @("HD",   225.61,                1.6994, 80.33750188338104)
)

$r = 2
foreach ($row in $rows) {
    $watchlist.Cells.Item($r, 1).Value = $row[0]
    $watchlist.Cells.Item($r, 2).Value = $row[1]
    $watchlist.Cells.Item($r, 3).Value = $row[2]
    $watchlist.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. stocks — refreshed current_price / value / performance / current_rsi,
#    and RTX's purch_price + performance are now known (were "NA")
# ---------------------------------------------------------------------------
$stocks = $wb.Worksheets.Item("stocks")

$stocks.Cells.Item(2, 3).Value = 96.84999999999999
$stocks.Cells.Item(2, 5).Value = 968.5
$stocks.Cells.Item(2, 6).Value = 0.8329
$stocks.Cells.Item(2, 7).Value = 31.12028158757131

$stocks.Cells.Item(3, 3).Value = 42.02
$stocks.Cells.Item(3, 5).Value = 1008.48
$stocks.Cells.Item(3, 6).Value = 0.8641
$stocks.Cells.Item(3, 7).Value = 25.58558558558559

$stocks.Cells.Item(4, 3).Value = 108.91
$stocks.Cells.Item(4, 5).Value = 762.37
$stocks.Cells.Item(4, 6).Value = -1.1706
$stocks.Cells.Item(4, 7).Value = 30.06517427033153

$stocks.Cells.Item(5, 2).Value = 59.04999999999999
$stocks.Cells.Item(5, 3).Value = 59
$stocks.Cells.Item(5, 5).Value = 708
$stocks.Cells.Item(5, 6).Value = -0.0847
$stocks.Cells.Item(5, 7).Value = 29.25180519280995

# ---------------------------------------------------------------------------
# 3. portfolio — stocks / total values
# ---------------------------------------------------------------------------
$portfolio = $wb.Worksheets.Item("portfolio")
$portfolio.Cells.Item(3, 2).Value = 3447.35
$portfolio.Cells.Item(4, 2).Value = 10007.01

# ---------------------------------------------------------------------------
# 4. summary — overwrite the last two log rows with the refreshed totals,
#    and drop the now-redundant trailing row
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("summary")

$summary.Cells.Item(5, 1).Value = "05/05/2020 16:35:30"
$summary.Cells.Item(5, 2).Value = 6559.66
$summary.Cells.Item(5, 3).Value = 3447.35
$summary.Cells.Item(5, 4).Value = 10007.01

$summary.Cells.Item(6, 1).Value = "06/05/2020 06:20:52"
$summary.Cells.Item(6, 2).Value = 6559.66
$summary.Cells.Item(6, 3).Value = 3447.35
$summary.Cells.Item(6, 4).Value = 10007.01

$summary.Rows.Item(7).Delete()
